$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string (e.g. "211.92")
# need to be forced to Text format first, otherwise Excel would silently
# coerce them into floating point numbers (losing the exact text, e.g.
# trailing zeros like "230.10" -> 230.1) - matching the source data which
# stores every Price/Volume cell as literal text.
$textForceCells = @('D5', 'D8', 'D9', 'D11', 'D18', 'D19', 'D24', 'D26', 'D27', 'D41', 'D43', 'D44', 'D45', 'D47', 'D48', 'D51')
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values for this refresh.
$ws.Range('D2').Value2 = '28.311.14'
$ws.Range('E2').Value2 = '  -0.74%  '
$ws.Range('D3').Value2 = '1.574.58'
$ws.Range('E3').Value2 = '  -0.16%  '
$ws.Range('E4').Value2 = '  +0.11%  '
$ws.Range('D5').Value2 = '211.92'
$ws.Range('E5').Value2 = '  -0.33%  '
$ws.Range('E6').Value2 = '  -0.67%  '
$ws.Range('E7').Value2 = '  +0.15%  '
$ws.Range('D8').Value2 = '44.49'
$ws.Range('E8').Value2 = '  -3.82%  '
$ws.Range('D9').Value2 = '23.80'
$ws.Range('E9').Value2 = '  -0.98%  '
$ws.Range('E10').Value2 = '  -0.95%  '
$ws.Range('D11').Value2 = '0.0587'
$ws.Range('E11').Value2 = '  -1.01%  '
$ws.Range('E12').Value2 = '  +1.46%  '
$ws.Range('D13').Value2 = '1.800.46'
$ws.Range('E13').Value2 = '  -0.07%  '
$ws.Range('D14').Value2 = '1.582.16'
$ws.Range('E14').Value2 = '  +0.22%  '
$ws.Range('E15').Value2 = '  -0.63%  '
$ws.Range('D17').Value2 = '28.343.60'
$ws.Range('E17').Value2 = '  -0.60%  '
$ws.Range('D18').Value2 = '61.58'
$ws.Range('E18').Value2 = '  -1.43%  '
$ws.Range('D19').Value2 = '230.10'
$ws.Range('E19').Value2 = '  +0.24%  '
$ws.Range('E20').Value2 = '  +0.01%  '
$ws.Range('D21').Value2 = '0.0₃0684'
$ws.Range('E21').Value2 = '  -1.48%  '
$ws.Range('E22').Value2 = '  +0.10%  '
$ws.Range('E23').Value2 = '  +0.35%  '
$ws.Range('D24').Value2 = '9.04'
$ws.Range('E24').Value2 = '  -1.45%  '
$ws.Range('E25').Value2 = '  +1.21%  '
$ws.Range('D26').Value2 = '151.67'
$ws.Range('E26').Value2 = '  +0.10%  '
$ws.Range('D27').Value2 = '14.93'
$ws.Range('E27').Value2 = '  -0.72%  '
$ws.Range('E28').Value2 = '  -1.78%  '
$ws.Range('E29').Value2 = '  -1.52%  '
$ws.Range('E30').Value2 = '  +0.11%  '
$ws.Range('E31').Value2 = '  +3.18%  '
$ws.Range('E32').Value2 = '  -3.74%  '
$ws.Range('E33').Value2 = '  -0.48%  '
$ws.Range('E34').Value2 = '  -2.18%  '
$ws.Range('D35').Value2 = '1.383.99'
$ws.Range('E35').Value2 = '  -0.89%  '
$ws.Range('E36').Value2 = '  +5.74%  '
$ws.Range('E37').Value2 = '  -3.24%  '
$ws.Range('E38').Value2 = '  +0.03%  '
$ws.Range('E39').Value2 = '  +2.86%  '
$ws.Range('E40').Value2 = '  -1.82%  '
$ws.Range('D41').Value2 = '0.518'
$ws.Range('E41').Value2 = '  -3.09%  '
$ws.Range('E42').Value2 = '  +0.09%  '
$ws.Range('D43').Value2 = '1.88'
$ws.Range('E43').Value2 = '  +1.61%  '
$ws.Range('D44').Value2 = '0.786'
$ws.Range('E44').Value2 = '  -1.25%  '
$ws.Range('D45').Value2 = '0.0463'
$ws.Range('E45').Value2 = '  +0.34%  '
$ws.Range('E46').Value2 = '  -4.33%  '
$ws.Range('D47').Value2 = '0.925'
$ws.Range('E47').Value2 = '  -5.60%  '
$ws.Range('D48').Value2 = '62.30'
$ws.Range('E48').Value2 = '  -0.61%  '
$ws.Range('D49').Value2 = '1.712.10'
$ws.Range('E49').Value2 = '  -0.04%  '
$ws.Range('E50').Value2 = '  +0.62%  '
$ws.Range('D51').Value2 = '85.51'
$ws.Range('E51').Value2 = '  -0.18%  '

# Restore the default (Normal) style on the cells we temporarily marked as
# Text, so formatting matches the original workbook.
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
